$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Replace the leading "Video 9" (runs: "Video", " ", "9") with a single
# "Video09a" run, leaving the rest of the title ("- Validity and reliability")
# untouched.
$sub = $tr.Characters(1, 7)
$sub.Text = "Video09a"
